$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -1
$ws.Range("B1").Value = 3.647341012954712
$ws.Range("C1").Value = 1.960239887237549
$ws.Range("D1").Value = 1.536244750022888
$ws.Range("E1").Value = 1.396791100502014
